$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.770.35'
$ws.Range("E2").Value = '  -1.01%  '
$ws.Range("D3").Value = '2.355.65'
$ws.Range("E3").Value = '  -1.22%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.65'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '107.95'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.63%  '
$ws.Range("E7").Value = '  -1.54%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.621'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.24%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.12%  '
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.47'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.998'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.29%  '
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.98'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.23%  '
$ws.Range("D16").Value = '2.712.18'
$ws.Range("E16").Value = '  -1.06%  '
$ws.Range("D17").Value = '2.288.76'
$ws.Range("E17").Value = '  -3.93%  '
$ws.Range("D18").Value = '42.762.81'
$ws.Range("E18").Value = '  -0.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.45%  '
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '76.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.47%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '257.59'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.12%  '
$ws.Range("E24").Value = '  -4.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.91%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.42'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.14%  '
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("E29").Value = '  +1.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '175.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.27%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '36.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0891'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.57%  '
$ws.Range("E34").Value = '  -8.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.127'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +18.50%  '
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("E37").Value = '  -5.28%  '
$ws.Range("E38").Value = '  -1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.82'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.67'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.239'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '71.77'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.39%  '
$ws.Range("E43").Value = '  -7.48%  '
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '113.96'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.01'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.52%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.50'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.58'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.79%  '

# Row 51: TrustWalletToken -> Cronos
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.100'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.77%  '
